$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new resale-number row for 2025-02-25 08:48 at row 88
$row = 88

# Columns A and D contain text that Excel would otherwise auto-coerce
# (A looks like a date, D looks like a leading-zero number), so force
# them to Text format before writing, then restore the default "Normal"
# style so no stray formatting is left behind on the new row.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-25"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "08:48:06"

$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "08"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 130623
$ws.Cells.Item($row, 6).Value = 141895
$ws.Cells.Item($row, 7).Value = 172653
$ws.Cells.Item($row, 8).Value = 158951
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146461
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193858
$ws.Cells.Item($row, 14).Value = 115474
$ws.Cells.Item($row, 15).Value = 46544
$ws.Cells.Item($row, 16).Value = 29405
$ws.Cells.Item($row, 17).Value = 68575
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47986
$ws.Cells.Item($row, 20).Value = -1
